$wb = $excel.ActiveWorkbook

$wsBugs = $wb.Worksheets.Item("Bugs and small tasks")
$wsStories = $wb.Worksheets.Item("Stories")
$wsContent = $wb.Worksheets.Item("Content")
$wsUI = $wb.Worksheets.Item("UI")

# --- Content sheet (sheet4): new row 38 ---
$wsContent.Range("A38").Value = "End of game cut-scene"
$wsContent.Range("F38").Value = "4-post bed, or bench"
$wsContent.Rows.Item(38).RowHeight = 30

# --- Bugs and small tasks (sheet1): new row 2, A2 ---
$wsBugs.Range("A2").Value = "Gloves glow when an enemy is around or when a power up is used"

# --- Stories sheet (sheet3): new rows ---
$wsStories.Range("A40").Value = "Hero stats"
$wsStories.Range("F40").Value = "Number of times revived, time to complete, etc"
$wsStories.Rows.Item(40).RowHeight = 30

$wsStories.Range("A25").Value = "Connect to FB reminder popup"

$wsStories.Range("A41").Value = "Video replay"
$wsStories.Range("F41").Value = "Everyplay, ReplayKit for IOS"

$wsStories.Range("A42").Value = "More hero model options"

$wsStories.Range("A43").Value = "Hero's look customizable"

$wsStories.Range("A44").Value = "Secret lair Easter egg"

# --- Bugs and small tasks (sheet1): F2 ---
$wsBugs.Range("F2").Value = "Runes light up when enemies are nearby, or you are in heroic or legendary modes or you have purchased the coin doubler, or you are using a power-up"
$wsBugs.Rows.Item(2).RowHeight = 75

$wsStories.Range("A45").Value = "Ziplining"

$wsStories.Range("A46").Value = "Endless, distance-based mode that gets unlocked after main quest is completed"
$wsStories.Rows.Item(46).RowHeight = 30

$wsStories.Range("A47").Value = "Collect for Runes to be able to activate Cullis Gate"
$wsStories.Rows.Item(47).RowHeight = 30

$wsStories.Range("A48").Value = "Loading Menu tips"

$wsStories.Range("A49").Value = "Fortune Teller gipsy caravan"
$wsStories.Range("F49").Value = "Subsribe to newsletter, get quirky horoscope"
$wsStories.Rows.Item(49).RowHeight = 30

$wsStories.Range("A50").Value = "Keep troll as pursuier or simply use as an enemy?"
$wsStories.Rows.Item(50).RowHeight = 30

# --- Row height tweaks on existing rows ---
$wsStories.Rows.Item(13).RowHeight = 30
$wsStories.Rows.Item(18).RowHeight = 30

$wsUI.Rows.Item(2).RowHeight = 75
$wsUI.Rows.Item(3).RowHeight = 75
$wsUI.Rows.Item(4).RowHeight = 30

# --- Selections / view state ---
$wsBugs.Range("G6").Select()
$wsStories.Range("A53").Select()
$wsContent.Range("D24").Select()
